$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("g11.1")

$ws.Range("A9").Value = 2023
$ws.Range("B9").Value = 3.118144130554446
$ws.Range("C9").Value = 6.515074339641291
